# Add a new "break_on_off" column (L) to Sheet1 with header + 72 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Values for L2:L73 (row 1 is header). Rows 19, 37 and 54 are 1, all others 0.
$values = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$ws.Cells.Item(1, 12).Value = "break_on_off"

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $values[$i]
}

# Reset the view: scroll back to top-left and select L1:L73
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("L1:L73").Select()
